# Remove the "Problemi" section (heading + its two body paragraphs) and the
# blank paragraph that preceded it, which followed the "Link alla home"
# paragraph ending in "...?page=index". The collapsed "_GoBack" bookmark
# that originally sat inside the removed text is recreated at the end of
# the surviving paragraph.

$d = $word.ActiveDocument

# Locate the "Link alla home" Heading 2 paragraph; the paragraph right
# after it is the one we keep (it holds the "http://server/ index.php
# ?page=index" text and must remain the last paragraph of the document).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Link alla home") {
        $anchorIndex = $i
    }
}

$keepParaIndex = $anchorIndex + 1
$keepPara = $d.Paragraphs.Item($keepParaIndex)

# A Range collapsed exactly at "end of paragraph, just before the pilcrow"
# confuses Bookmarks.Add in this host (it silently snaps to some unrelated
# paragraph). Work around it: insert a throwaway character into the
# (otherwise blank) paragraph that immediately follows the "keep"
# paragraph, wrap a bookmark around that single character (a non-empty
# Range works fine), then delete the "keep" paragraph's own paragraph
# mark so the marked character is pulled into the "keep" paragraph's text
# -- this leaves the "keep" run's rsid untouched, unlike inserting
# straight into it.
$nextPara = $d.Paragraphs.Item($keepParaIndex + 1)
$insertPoint = $d.Range($keepPara.Range.End, $keepPara.Range.End)
$insertPoint.InsertBefore("X")

$nextPara = $d.Paragraphs.Item($keepParaIndex + 1)
$markerRange = $d.Range($nextPara.Range.Start, $nextPara.Range.Start + 1)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $markerRange)

$keepPara = $d.Paragraphs.Item($keepParaIndex)
$mergeRange = $d.Range($keepPara.Range.End - 1, $keepPara.Range.End)
$mergeRange.Delete()

# Delete everything from the start of the (now following) blank paragraph
# through to the end of the document body content: this removes the
# original blank paragraph's leftovers, the "Problemi" heading and the
# two paragraphs describing the install problem / alternate link.
$firstToRemove = $d.Paragraphs.Item($keepParaIndex + 1)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$deleteRange = $d.Range($firstToRemove.Range.Start, $lastPara.Range.End)
$deleteRange.Delete()

# Clear the throwaway marker character via the bookmark's own Range --
# this leaves "_GoBack" correctly collapsed right before the paragraph
# mark of the (now last) surviving paragraph.
$d.Bookmarks.Item("_GoBack").Range.Text = ""

Write-Output ("Done. ParaCount=" + $d.Paragraphs.Count)
